$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Actualiza base de datos EC: el periodo en mora pasa de 2507 a 2508
# para todos los trabajadores relacionados (columna "Periodo Mora", filas 16-20)
$ws.Range("E16:E20").Value = "2508"
